$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "NSG-1522"
$ws.Range("B2").Value = "1522-NSG"
$ws.Range("C2").Value = "HERO NAME WALA STIKER"
$ws.Range("D2").Value = "SPL I3S"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 250
$ws.Range("G2").Value = "MS"
